# Correçao de bugs em CP e Reatores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Resistores de descarga"

# Remove extra data rows (old rows 5-10), keeping only rows 2-4
$ws.Rows("5:10").Delete()

# Remove the "Número de resitores em paralelo" column (old column C),
# shifting everything after it one column to the left
$ws.Columns(3).Delete()

# Update header row (row 1) - after the column delete, headers shifted:
# old D1 "Resistor usado" is now C1, etc. Overwrite each with new text.
$ws.Range("C1").Value = "Resistor usado"
$ws.Range("D1").Value = "Restência total"
$ws.Range("E1").Value = "Max WATTS com 110% Vn"
$ws.Range("F1").Value = "110% da Resistência total"
$ws.Range("G1").Value = "90% da Resistência total"
$ws.Range("H1").Value = "Resistor máximo permitido"

# Update data rows 2-4 with the new computed values
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 0.96
$ws.Range("D2").Value = 5.76
$ws.Range("E2").Value = 14.01088011695907
$ws.Range("F2").Value = 6.336
$ws.Range("G2").Value = 5.184
$ws.Range("H2").Value = 6.299725638800916

$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 0.55
$ws.Range("D3").Value = 3.850000000000001
$ws.Range("E3").Value = 20.96173233082706
$ws.Range("F3").Value = 4.235000000000001
$ws.Range("G3").Value = 3.465
$ws.Range("H3").Value = 6.299725638800916

$ws.Range("B4").Value = 11
$ws.Range("C4").Value = 0.23
$ws.Range("D4").Value = 2.53
$ws.Range("E4").Value = 31.89828832951946
$ws.Range("F4").Value = 2.783
$ws.Range("G4").Value = 2.277
$ws.Range("H4").Value = 6.299725638800916
